$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for Il6-Il6st LR-pair table (rows 2-10)
$ws.Range("G2").Value = 0.3056103333333333
$ws.Range("H2").Value = 0.916831
$ws.Range("I2").Value = 0.01726097181671177
$ws.Range("J2").Value = 0.01726097181671177
$ws.Range("M2").Value = 20.854426
$ws.Range("N2").Value = 62.563278
$ws.Range("O2").Value = 0.1507164072139519
$ws.Range("P2").Value = 0.1507164072139519
$ws.Range("Q2").Value = 6.373328081335333
$ws.Range("R2").Value = 57.35995273201799
$ws.Range("S2").Value = 0.002601511657236079
$ws.Range("T2").Value = 0.002601511657236079
$ws.Range("G3").Value = 0.3056103333333333
$ws.Range("H3").Value = 0.916831
$ws.Range("I3").Value = 0.01726097181671177
$ws.Range("J3").Value = 0.01726097181671177
$ws.Range("O3").Value = 0.6862909728343718
$ws.Range("P3").Value = 0.6862909728343718
$ws.Range("Q3").Value = 29.02111064074878
$ws.Range("R3").Value = 261.189995766739
$ws.Range("S3").Value = 0.0118460491401578
$ws.Range("T3").Value = 0.0118460491401578
$ws.Range("G4").Value = 0.3056103333333333
$ws.Range("H4").Value = 0.916831
$ws.Range("I4").Value = 0.01726097181671177
$ws.Range("J4").Value = 0.01726097181671177
$ws.Range("N4").Value = 67.65920700000001
$ws.Range("O4").Value = 0.1629926199516763
$ws.Range("P4").Value = 0.1629926199516763
$ws.Range("Q4").Value = 6.892450934779667
$ws.Range("R4").Value = 62.032058413017
$ws.Range("S4").Value = 0.002813411019317897
$ws.Range("T4").Value = 0.002813411019317897
$ws.Range("I5").Value = 0.8433360339088308
$ws.Range("J5").Value = 0.8433360339088307
$ws.Range("M5").Value = 20.854426
$ws.Range("N5").Value = 62.563278
$ws.Range("O5").Value = 0.1507164072139519
$ws.Range("P5").Value = 0.1507164072139519
$ws.Range("Q5").Value = 311.3878687704753
$ws.Range("R5").Value = 2802.490818934278
$ws.Range("S5").Value = 0.1271045771048025
$ws.Range("T5").Value = 0.1271045771048025
$ws.Range("I6").Value = 0.8433360339088308
$ws.Range("J6").Value = 0.8433360339088307
$ws.Range("O6").Value = 0.6862909728343718
$ws.Range("P6").Value = 0.6862909728343718
$ws.Range("S6").Value = 0.5787739071375723
$ws.Range("T6").Value = 0.5787739071375722
$ws.Range("I7").Value = 0.8433360339088308
$ws.Range("J7").Value = 0.8433360339088307
$ws.Range("N7").Value = 67.65920700000001
$ws.Range("O7").Value = 0.1629926199516763
$ws.Range("P7").Value = 0.1629926199516763
$ws.Range("Q7").Value = 336.7511572911897
$ws.Range("S7").Value = 0.137457549666456
$ws.Range("T7").Value = 0.137457549666456
$ws.Range("I8").Value = 0.1394029942744574
$ws.Range("J8").Value = 0.1394029942744574
$ws.Range("M8").Value = 20.854426
$ws.Range("N8").Value = 62.563278
$ws.Range("O8").Value = 0.1507164072139519
$ws.Range("P8").Value = 0.1507164072139519
$ws.Range("Q8").Value = 51.472247765994
$ws.Range("R8").Value = 463.250229893946
$ws.Range("S8").Value = 0.02101031845191333
$ws.Range("T8").Value = 0.02101031845191333
$ws.Range("I9").Value = 0.1394029942744574
$ws.Range("J9").Value = 0.1394029942744574
$ws.Range("O9").Value = 0.6862909728343718
$ws.Range("P9").Value = 0.6862909728343718
$ws.Range("S9").Value = 0.09567101655664172
$ws.Range("T9").Value = 0.09567101655664172
$ws.Range("I10").Value = 0.1394029942744574
$ws.Range("J10").Value = 0.1394029942744574
$ws.Range("N10").Value = 67.65920700000001
$ws.Range("O10").Value = 0.1629926199516763
$ws.Range("P10").Value = 0.1629926199516763
$ws.Range("Q10").Value = 55.66478576066101
$ws.Range("R10").Value = 500.983071845949
$ws.Range("S10").Value = 0.02272165926590233
$ws.Range("T10").Value = 0.02272165926590233
